# Auto-generated edit script: adds a new 'handoff' row (row 3) to the
# Overview / zh-cn / de-de sheets, resizes their tables, adds the two new
# hyperlinks, and widens the affected status columns.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet (sheet1) ----
$wsOverview = $wb.Worksheets.Item(1)

$wsOverview.Range('A3').Formula = "'" + '3a80ca73-81d2-4e98-9796-591b3bf3be15ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$wsOverview.Range('C3').Formula = "'" + '.md'
$wsOverview.Range('D3').Formula = "'"
$wsOverview.Range('E3').Formula = "'" + 'Ready for handoff'
$wsOverview.Range('F3').Formula = "'" + 'Ready for handoff'
$wsOverview.Range('G3').Formula = '2016-08-22 20:27:20'
$wsOverview.Hyperlinks.Add($wsOverview.Range('B3'), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/78d22917b60c1478978e9d9cf7b47b59c30edfb7/e2e/3a80ca73-81d2-4e98-9796-591b3bf3be15ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md', '', '', 'e2e\3a80ca73-81d2-4e98-9796-591b3bf3be15ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md')

# resize table + autofilter to new extent
$wsOverview.ListObjects.Item(1).Resize($wsOverview.Range('A1:G3'))

# widen zh-cn / de-de status columns (E, F) to match new content width
$wsOverview.Columns.Item(5).ColumnWidth = 16.25
$wsOverview.Columns.Item(6).ColumnWidth = 16.25

# ---- zh-cn sheet (sheet2) ----
$wsZhCn = $wb.Worksheets.Item(2)

$wsZhCn.Range('B3').Formula = '.md'
$wsZhCn.Range('C3').Formula = 'Ready for handoff'
$wsZhCn.Range('D3').Formula = 'e2e'
$wsZhCn.Range('E3').Formula = 'ht'
$wsZhCn.Range('F3').Formula = "'" + 'False'
$wsZhCn.Range('G3').Formula = '3a80ca73-81d2-4e98-9796-591b3bf3be15oooooooooooooooooooooooooooooooooooooooo.f9e01823d1840005a7b4d64e594ea6923f490efc.zh-cn.xlf'
$wsZhCn.Range('H3').Formula = '2016-08-22 20:27:11'
$wsZhCn.Range('I3').Formula = "'"
$wsZhCn.Range('J3').Formula = "'"
$wsZhCn.Range('K3').Formula = '0001-01-01 00:00:00'
$wsZhCn.Range('L3').Formula = "'"
$wsZhCn.Range('M3').Formula = "'" + 'True'
$wsZhCn.Range('N3').Formula = "'"
$wsZhCn.Range('O3').Formula = "'" + 'False'
$wsZhCn.Range('P3').Formula = "'"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range('A3'), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/78d22917b60c1478978e9d9cf7b47b59c30edfb7/e2e/3a80ca73-81d2-4e98-9796-591b3bf3be15ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md', '', '', '3a80ca73-81d2-4e98-9796-591b3bf3be15ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md')

$wsZhCn.ListObjects.Item(1).Resize($wsZhCn.Range('A1:P3'))
$wsZhCn.Columns.Item(3).ColumnWidth = 16.25

# ---- de-de sheet (sheet3) ----
$wsDeDe = $wb.Worksheets.Item(3)

$wsDeDe.Range('B3').Formula = '.md'
$wsDeDe.Range('C3').Formula = 'Ready for handoff'
$wsDeDe.Range('D3').Formula = 'e2e'
$wsDeDe.Range('E3').Formula = 'ht'
$wsDeDe.Range('F3').Formula = "'" + 'False'
$wsDeDe.Range('G3').Formula = '3a80ca73-81d2-4e98-9796-591b3bf3be15oooooooooooooooooooooooooooooooooooooooo.f9e01823d1840005a7b4d64e594ea6923f490efc.de-de.xlf'
$wsDeDe.Range('H3').Formula = '2016-08-22 20:27:20'
$wsDeDe.Range('I3').Formula = "'"
$wsDeDe.Range('J3').Formula = "'"
$wsDeDe.Range('K3').Formula = '0001-01-01 00:00:00'
$wsDeDe.Range('L3').Formula = "'"
$wsDeDe.Range('M3').Formula = "'" + 'True'
$wsDeDe.Range('N3').Formula = "'"
$wsDeDe.Range('O3').Formula = "'" + 'False'
$wsDeDe.Range('P3').Formula = "'"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range('A3'), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/78d22917b60c1478978e9d9cf7b47b59c30edfb7/e2e/3a80ca73-81d2-4e98-9796-591b3bf3be15ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md', '', '', '3a80ca73-81d2-4e98-9796-591b3bf3be15ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md')

$wsDeDe.ListObjects.Item(1).Resize($wsDeDe.Range('A1:P3'))
$wsDeDe.Columns.Item(3).ColumnWidth = 16.25

